$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: fix(publipostage): Try to solve Excel emoji problem
# Replace the "statut" marker emoji used throughout column A with new
# plain-text / emoji markers:
#   📘 -> ⚠️
#   📕 -> -3
#   📙 -> +3
#   📗 -> ✅

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $v = $cell.Value2

    if ($v -eq "📘") {
        $cell.Value = "⚠️"
    } elseif ($v -eq "📕") {
        # "-3" looks numeric to Excel's auto-detection, so force text
        # storage the same way typing an apostrophe-prefixed / text
        # formatted value in the UI would.
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
    } elseif ($v -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
    } elseif ($v -eq "📗") {
        $cell.Value = "✅"
    }
}

Write-Host "done updating statut column"
